# Automatische test-sync: 2025-07-31 21:46:50
# Adds the new "Testmail #11" row to the Logs sheet, rolls the matching
# category count into the Dashboard sheet, and extends the conditional
# formatting + chart series ranges so they keep covering the new row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 13 with the new mail log entry
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A13").Value = "Mijn retour is nog steeds niet verwerkt."
$logs.Range("B13").Value = "mailmind.test@zohomail.eu"
$logs.Range("C13").Value = "Testmail #11: Mijn retour is nog steeds niet verwerkt."
$logs.Range("D13").Value = "Retour / Terugbetaling"
$logs.Range("E13").Value = "Beste klant,`nBedankt voor je bericht. We begrijpen dat het vervelend is dat je retourzending nog niet verwerkt is. Om je verder te kunnen helpen, ontvangen we graag wat aanvullende informatie zoals het ordernummer van de retourzending. Zodra we deze gegevens hebben, zullen we direct voor je aan de slag gaan om het probleem op te lossen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F13").Value = "2025-07-31 21:45:55"
$logs.Range("G13").Value = "Ja"
$logs.Range("H13").Value = "Nee"
$logs.Range("I13").Value = "Ja"
$logs.Range("J13").Value = "Nee"

# The multi-line text in E13 otherwise leaves the row with an explicit
# custom height; AutoFit puts it back to the sheet's implicit default,
# matching the rest of the rows (and the target row 13).
$logs.Rows.Item(13).AutoFit()

# Extend conditional formatting ranges from row 12 to row 13
$logs.Range("D2:D12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D13"))
$logs.Range("G2:G12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G13"))
$logs.Range("H2:H12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H13"))
$logs.Range("I2:I12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I13"))
$logs.Range("J2:J12").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J13"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: add the new category count row 6
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Retour / Terugbetaling"
$dash.Range("B6").Value = 1

# ---------------------------------------------------------------------
# 3. Chart: extend the category/value series references to include row 6
# ---------------------------------------------------------------------
$chartObj = $dash.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$6"
$series.Values = "='Dashboard'!`$B`$2:`$B`$6"
